$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The codeforiati:group-code and codeforiati:group-name columns (C and D)
# have swapped places: what used to be column D (group-code) is now
# column C, and what used to be column C (group-name) is now column D.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}
